$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'316.16"
$ws.Range("E2").Value = "'3.51%"
$ws.Range("D3").Value = "'39.43"
$ws.Range("E3").Value = "'2.92%"
$ws.Range("D4").Value = "'5.119"
$ws.Range("E4").Value = "'0.52%"
$ws.Range("D5").Value = "'0.08203"
$ws.Range("E5").Value = "'1.80%"
$ws.Range("D6").Value = "'2.037"
$ws.Range("E6").Value = "'4.73%"
$ws.Range("E7").Value = "'3.89%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9332"
$ws.Range("E8").Value = "'0.35%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1410"
$ws.Range("E9").Value = "'-3.37%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1991"
$ws.Range("E10").Value = "'3.73%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09114"
$ws.Range("E11").Value = "'0.37%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03529"
$ws.Range("E12").Value = "'0.35%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09815"
$ws.Range("E13").Value = "'0.27%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001396"
$ws.Range("E14").Value = "'0.17%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006299"
$ws.Range("E15").Value = "'3.26%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.660"
$ws.Range("E16").Value = "'-1.71%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.280"
$ws.Range("E17").Value = "'2.27%"
$ws.Range("D18").Value = "'3.295"
$ws.Range("E18").Value = "'-4.47%"
$ws.Range("D19").Value = "'0.3462"
$ws.Range("E19").Value = "'0.01%"
$ws.Range("E20").Value = "'-0.54%"
$ws.Range("D21").Value = "'4.894"
$ws.Range("E21").Value = "'1.75%"
$ws.Range("D22").Value = "'0.2452"
$ws.Range("E22").Value = "'1.91%"
$ws.Range("D23").Value = "'0.04335"
$ws.Range("E23").Value = "'-0.55%"
$ws.Range("D24").Value = "'0.001225"
$ws.Range("E24").Value = "'-0.52%"
$ws.Range("D25").Value = "'0.004780"
$ws.Range("E25").Value = "'12.42%"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'-0.19%"
$ws.Range("D27").Value = "'0.0004002"
$ws.Range("E27").Value = "'-10.01%"
$ws.Range("D39").Value = "'0.02229"
$ws.Range("E39").Value = "'7.97%"
$ws.Range("D40").Value = "'0.05246"
$ws.Range("E40").Value = "'4.16%"
$ws.Range("D41").Value = "'0.007520"
$ws.Range("E41").Value = "'0.71%"
$ws.Range("D42").Value = "'0.009776"
$ws.Range("E42").Value = "'-3.36%"
$ws.Range("D43").Value = "'0.1378"
$ws.Range("E43").Value = "'2.33%"
$ws.Range("D44").Value = "'0.002149"
$ws.Range("E44").Value = "'0.27%"
$ws.Range("D45").Value = "'0.009502"
$ws.Range("E45").Value = "'6.64%"
$ws.Range("D46").Value = "'0.00006456"
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("E48").Value = "'-25.07%"
$ws.Range("D49").Value = "'0.002768"
$ws.Range("E49").Value = "'-1.91%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.19%"
